$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.968.30"
$ws.Range("E2").Value = "  -3.26%  "

$ws.Range("D3").Value = "3.182.59"
$ws.Range("E3").Value = "  -3.43%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.60%  "

$ws.Range("E7").Value = "  -5.74%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").Value = "3.180.65"
$ws.Range("E9").Value = "  -3.39%  "

$ws.Range("E10").Value = "  -3.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.46%  "

$ws.Range("E12").Value = "  -4.58%  "

$ws.Range("D13").Value = "3.737.24"
$ws.Range("E13").Value = "  -3.35%  "

$ws.Range("E14").Value = "  -1.95%  "

$ws.Range("D15").Value = "64.086.36"
$ws.Range("E15").Value = "  -3.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.92%  "

$ws.Range("E17").Value = "  -2.47%  "

$ws.Range("D18").Value = "3.193.42"
$ws.Range("E18").Value = "  -2.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "415.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.28%  "

$ws.Range("E25").Value = "  +3.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.494"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.46%  "

$ws.Range("E27").Value = "  -3.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.84%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.37%  "

$ws.Range("E35").Value = "  -4.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "155.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.52%  "

$ws.Range("D38").Value = "2.746.50"
$ws.Range("E38").Value = "  -1.45%  "

$ws.Range("E39").Value = "  -4.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.718"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0629"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.07%  "

$ws.Range("E47").Value = "  -1.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "298.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0996"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.36%  "

$ws.Range("E51").Value = "  +0.00%  "
